$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 600
$ws.Range("I31").Value = 600
$ws.Range("K31").Value = 1800
$ws.Range("M31").Value = -1570
$ws.Range("H33").Value = 556.875
$ws.Range("I33").Value = 639.2308
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 639.2308
$ws.Range("L33").Value = 200
$ws.Range("M33").Value = -410.2308
$ws.Range("N33").Value = -658
$ws.Range("H64").Value = 3077.2727
$ws.Range("I64").Value = 3065
$ws.Range("J64").Value = 3200
$ws.Range("K64").Value = 3065
$ws.Range("L64").Value = 3200
$ws.Range("M64").Value = -2817
$ws.Range("N64").Value = -3696
$ws.Range("H67").Value = 3077.2727
$ws.Range("I67").Value = 3065
$ws.Range("J67").Value = 3200
$ws.Range("K67").Value = 3065
$ws.Range("L67").Value = 3200
$ws.Range("M67").Value = -2207
$ws.Range("N67").Value = -4916
$ws.Range("H74").Value = 3738.5
$ws.Range("I74").Value = 3075
$ws.Range("J74").Value = 3928.0715
$ws.Range("K74").Value = 3075
$ws.Range("L74").Value = 3928.0715
$ws.Range("M74").Value = -2139
$ws.Range("N74").Value = -5800.0715
$ws.Range("H76").Value = 187857.58
$ws.Range("I76").Value = 169167.17
$ws.Range("J76").Value = 300000
$ws.Range("K76").Value = 169167.17
$ws.Range("L76").Value = 300000
$ws.Range("M76").Value = -168852.17
$ws.Range("N76").Value = -300630
$ws.Range("H77").Value = 3738.5
$ws.Range("I77").Value = 3075
$ws.Range("J77").Value = 3928.0715
$ws.Range("K77").Value = 15375
$ws.Range("L77").Value = 19640.3575
$ws.Range("M77").Value = -10695
$ws.Range("N77").Value = -29000.3575
$ws.Range("H79").Value = 187857.58
$ws.Range("I79").Value = 169167.17
$ws.Range("J79").Value = 300000
$ws.Range("K79").Value = 169167.17
$ws.Range("L79").Value = 300000
$ws.Range("M79").Value = -168075.17
$ws.Range("N79").Value = -302184
$ws.Range("H137").Value = 1467.7812
$ws.Range("I137").Value = 1328.1666
$ws.Range("J137").Value = 1886.625
$ws.Range("K137").Value = 3984.4998
$ws.Range("L137").Value = 5659.875
$ws.Range("M137").Value = -1434.4998
$ws.Range("N137").Value = -10759.875

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1883.0834
$ws.Range("I2").Value = 1225
$ws.Range("J2").Value = 2212.125
$ws.Range("K2").Value = 1225
$ws.Range("L2").Value = 2212.125
$ws.Range("M2").Value = -1112
$ws.Range("N2").Value = -2438.125
$ws.Range("H61").Value = 2041.5
$ws.Range("I61").Value = 1199.5714
$ws.Range("J61").Value = 2577.2727
$ws.Range("K61").Value = 1199.5714
$ws.Range("L61").Value = 2577.2727
$ws.Range("M61").Value = -987.5714
$ws.Range("N61").Value = -3001.2727
$ws.Range("H63").Value = 1669750
$ws.Range("I63").Value = 1669750
$ws.Range("K63").Value = 1669750
$ws.Range("M63").Value = -1669064
$ws.Range("H66").Value = 1669750
$ws.Range("I66").Value = 1669750
$ws.Range("K66").Value = 8348750
$ws.Range("M66").Value = -8345318
$ws.Range("H74").Value = 786.80554
$ws.Range("I74").Value = 762.5357
$ws.Range("J74").Value = 871.75
$ws.Range("K74").Value = 762.5357
$ws.Range("L74").Value = 871.75
$ws.Range("M74").Value = 111.4643
$ws.Range("N74").Value = -2619.75
$ws.Range("H77").Value = 786.80554
$ws.Range("I77").Value = 762.5357
$ws.Range("J77").Value = 871.75
$ws.Range("K77").Value = 3812.6785
$ws.Range("L77").Value = 4358.75
$ws.Range("M77").Value = 555.3215
$ws.Range("N77").Value = -13094.75
$ws.Range("H116").Value = 1883.0834
$ws.Range("I116").Value = 1225
$ws.Range("J116").Value = 2212.125
$ws.Range("K116").Value = 1225
$ws.Range("L116").Value = 2212.125
$ws.Range("M116").Value = 1069
$ws.Range("N116").Value = -6800.125
$ws.Range("H132").Value = 4247.5264
$ws.Range("I132").Value = 6231.5
$ws.Range("K132").Value = 18694.5
$ws.Range("M132").Value = -16164.5
$ws.Range("H136").Value = 2041.5
$ws.Range("I136").Value = 1199.5714
$ws.Range("J136").Value = 2577.2727
$ws.Range("K136").Value = 3598.7142
$ws.Range("L136").Value = 7731.8181
$ws.Range("M136").Value = -1048.7142
$ws.Range("N136").Value = -12831.8181

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1883.0834
$ws.Range("I3").Value = 1225
$ws.Range("J3").Value = 2212.125
$ws.Range("K3").Value = 1225
$ws.Range("L3").Value = 2212.125
$ws.Range("M3").Value = -1111
$ws.Range("N3").Value = -2440.125
$ws.Range("H134").Value = 21166.809
$ws.Range("I134").Value = 30381.428
$ws.Range("J134").Value = 2195.5293
$ws.Range("K134").Value = 91144.284
$ws.Range("L134").Value = 6586.5879
$ws.Range("M134").Value = -88609.284
$ws.Range("N134").Value = -11656.5879

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7144411
$ws.Range("I31").Value = 1701.55
$ws.Range("J31").Value = 25001184
$ws.Range("K31").Value = 1701.55
$ws.Range("L31").Value = 25001184
$ws.Range("M31").Value = -1406.55
$ws.Range("N31").Value = -25001774
$ws.Range("H34").Value = 7144411
$ws.Range("I34").Value = 1701.55
$ws.Range("J34").Value = 25001184
$ws.Range("K34").Value = 1701.55
$ws.Range("L34").Value = 25001184
$ws.Range("M34").Value = -1499.55
$ws.Range("N34").Value = -25001588
$ws.Range("H58").Value = 1753.8462
$ws.Range("I58").Value = 2290
$ws.Range("J58").Value = 1418.75
$ws.Range("K58").Value = 2290
$ws.Range("L58").Value = 1418.75
$ws.Range("M58").Value = -2087
$ws.Range("N58").Value = -1824.75
$ws.Range("H132").Value = 2277.7568
$ws.Range("I132").Value = 2016.6316
$ws.Range("J132").Value = 2553.389
$ws.Range("K132").Value = 6049.8948
$ws.Range("L132").Value = 7660.167
$ws.Range("M132").Value = -3519.8948
$ws.Range("N132").Value = -12720.167
$ws.Range("H134").Value = 1292.4242
$ws.Range("I134").Value = 1246.8
$ws.Range("J134").Value = 1362.6154
$ws.Range("K134").Value = 3740.4
$ws.Range("L134").Value = 4087.8462
$ws.Range("M134").Value = -1205.4
$ws.Range("N134").Value = -9157.8462
$ws.Range("H136").Value = 1753.8462
$ws.Range("I136").Value = 2290
$ws.Range("J136").Value = 1418.75
$ws.Range("K136").Value = 6870
$ws.Range("L136").Value = 4256.25
$ws.Range("M136").Value = -4320
$ws.Range("N136").Value = -9356.25

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28337694
$ws.Range("I70").Value = 30004382
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 30004382
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -30004112
$ws.Range("N70").Value = -4540
$ws.Range("H73").Value = 28337694
$ws.Range("I73").Value = 30004382
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 30004382
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -30003446
$ws.Range("N73").Value = -5872
$ws.Range("H80").Value = 3760
$ws.Range("I80").Value = 3700
$ws.Range("J80").Value = 3768.5715
$ws.Range("K80").Value = 3700
$ws.Range("L80").Value = 3768.5715
$ws.Range("M80").Value = -2702
$ws.Range("N80").Value = -5764.5715
$ws.Range("H83").Value = 3760
$ws.Range("I83").Value = 3700
$ws.Range("J83").Value = 3768.5715
$ws.Range("K83").Value = 18500
$ws.Range("L83").Value = 18842.8575
$ws.Range("M83").Value = -13508
$ws.Range("N83").Value = -28826.8575
$ws.Range("H132").Value = 120866.18
$ws.Range("I132").Value = 224769.67
$ws.Range("K132").Value = 674309.01
$ws.Range("M132").Value = -671779.01

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 11479.315
$ws.Range("I132").Value = 20533.777
$ws.Range("J132").Value = 3330.3
$ws.Range("K132").Value = 61601.33099999999
$ws.Range("L132").Value = 9990.900000000001
$ws.Range("M132").Value = -59071.33099999999
$ws.Range("N132").Value = -15050.9
$ws.Range("H136").Value = 12988.889
$ws.Range("I136").Value = 25725
$ws.Range("K136").Value = 77175
$ws.Range("M136").Value = -74625

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1148.2368
$ws.Range("I132").Value = 1003.82855
$ws.Range("K132").Value = 3011.48565
$ws.Range("M132").Value = -481.4856499999996
$ws.Range("H136").Value = 6306.3477
$ws.Range("I136").Value = 7423.4736
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 22270.4208
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -19720.4208
$ws.Range("N136").Value = -8100
